$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: advance the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the price list values in column D (rows 31-38)
$ws.Range("D31").Value = 8213.654
$ws.Range("D32").Value = 7326.587
$ws.Range("D33").Value = 6986.684
$ws.Range("D34").Value = 11898.71
$ws.Range("D35").Value = 8878.953
$ws.Range("D36").Value = 8369.097
$ws.Range("D37").Value = 7842.661
$ws.Range("D38").Value = 12400.274
